$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (changed) date column C for rows 2-16 from 45207 to 45208
$ws.Range("C2:C16").Value = 45208
